$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.315.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.192.56'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.72%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.190.82'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.547'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.01'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.517'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.721.78'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.375.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.195.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '514.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.739'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.50%  '
$ws.Range('E28').Value = '  +4.20%  '
$ws.Range('E29').Value = '  +7.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +17.18%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.20%  '
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '507.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.86'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0900'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0424'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.88'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.20%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.302'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.05%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0670'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +16.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.916.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('E48').Value = '  +3.05%  '
$ws.Range('E50').Value = '  +5.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.82%  '
